$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$sub = $para1.Characters(6, 61)
$sub.Text = " 이번엔 암호화 프로그램을 만들었지만 다음 번에 기회가 된다면 충분한 시간을 가지고 파이썬 해킹 프로그래밍을 해보고 싶다"
